$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows above the current row 296, pushing the rest down.
$ws.Rows("296:297").Insert()

# New row 296
$ws.Cells.Item(296, 1).Value = 5
$ws.Cells.Item(296, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(296, 3).Value = "Maule"
$ws.Cells.Item(296, 4).Value = [DateTime]"2022-03-21"
$ws.Cells.Item(296, 5).Value = 7
$ws.Cells.Item(296, 6).Value = 100112043
$ws.Cells.Item(296, 7).Value = "Pepino ensalada"
$ws.Cells.Item(296, 8).Value = "Sin especificar"
$ws.Cells.Item(296, 9).Value = "Primera"
$ws.Cells.Item(296, 10).Value = 400
$ws.Cells.Item(296, 11).Value = 18000
$ws.Cells.Item(296, 12).Value = 18000
$ws.Cells.Item(296, 13).Value = 18000
$ws.Cells.Item(296, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(296, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(296, 16).Value = 300
$ws.Cells.Item(296, 17).Value = 60
$ws.Cells.Item(296, 18).Value = "Hortaliza"

# New row 297
$ws.Cells.Item(297, 1).Value = 5
$ws.Cells.Item(297, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(297, 3).Value = "Maule"
$ws.Cells.Item(297, 4).Value = [DateTime]"2022-03-21"
$ws.Cells.Item(297, 5).Value = 7
$ws.Cells.Item(297, 6).Value = 100112043
$ws.Cells.Item(297, 7).Value = "Pepino ensalada"
$ws.Cells.Item(297, 8).Value = "Sin especificar"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 350
$ws.Cells.Item(297, 11).Value = 19000
$ws.Cells.Item(297, 12).Value = 19000
$ws.Cells.Item(297, 13).Value = 19000
$ws.Cells.Item(297, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(297, 15).Value = "Región del Maule"
$ws.Cells.Item(297, 16).Value = 238
$ws.Cells.Item(297, 17).Value = 80
$ws.Cells.Item(297, 18).Value = "Hortaliza"
